$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "MSc DS 지원자의 착각"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/msc-ds-applicant-confusion/#utm_source=rss&utm_medium=rss&utm_campaign=msc-ds-applicant-confusion"

$ws.Range("D32").Value = "youtube Data API를 이용해 유튜브 댓글(라이브방송 포함) 수집"
$ws.Range("E32").Value = "https://dodonam.tistory.com/342"

$ws.Range("D50").Value = "DM21 (DeepMind 21)"
$ws.Range("E50").Value = "http://incredible.egloos.com/7531714"
